$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.2.0-ballot -> 2.1.0
$meta.Range("B3").Value = "2.1.0"

# Date: 2025-12-19T08:32:44+00:00 -> 2025-12-19T08:44:55+00:00
$meta.Range("B8").Value = "2025-12-19T08:44:55+00:00"

# Base Definition: drop the "|4.0.1" version suffix
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Binding Value Set: drop the "|2.2.0-ballot" version suffix
$elem.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/tddui/ValueSet/tddui-discriminator-vs"

# Column Z width shrinks (bestFit) because the new value is shorter than before.
# Target stored width is 59.12109375 characters; the COM ColumnWidth setter
# snaps to a 1/6-character pixel grid, so 58.35 is the closest input that
# lands on the nearest achievable grid point (59.16666...).
$elem.Columns.Item(26).ColumnWidth = 58.35
